$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: column G was "Website" (Mouser links); rename header to
#    "Mouser" and add a new column H "AliExpress" with AliExpress links.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Mouser"
$ws.Range("H1").Value = "AliExpress"

# ---------------------------------------------------------------------------
# 2. Populate the new AliExpress column (H) with plain-text URLs (no
#    hyperlink objects attached - same as how the Mouser column originally
#    held plain text in most rows).
# ---------------------------------------------------------------------------
$ws.Range("H2").Value  = "https://www.aliexpress.com/item/1068215251.html?spm=a2g0o.productlist.0.0.1e5ca144y6dBmp&algo_pvid=80ccebdb-5973-4931-b90a-dfee728a65f4&algo_expid=80ccebdb-5973-4931-b90a-dfee728a65f4-1&btsid=0ab50f4415858253351994678edc77&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H6").Value  = "https://www.aliexpress.com/item/33013972369.html?spm=a2g0o.productlist.0.0.25ae390djyGrti&algo_pvid=73949350-3a11-4e7e-9716-91a2fb3778bf&algo_expid=73949350-3a11-4e7e-9716-91a2fb3778bf-3&btsid=0ab6f81e15858255050474469e1747&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H7").Value  = "https://www.aliexpress.com/item/32966490820.html?spm=a2g0o.productlist.0.0.5cb5e375SYJ0QK&algo_pvid=4f258cba-cd74-4f42-ad43-205c371b1b84&algo_expid=4f258cba-cd74-4f42-ad43-205c371b1b84-0&btsid=0ab6f83915858255934824697e4e73&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H8").Value  = "https://www.aliexpress.com/item/32966490820.html?spm=a2g0o.productlist.0.0.5cb5e375SYJ0QK&algo_pvid=4f258cba-cd74-4f42-ad43-205c371b1b84&algo_expid=4f258cba-cd74-4f42-ad43-205c371b1b84-0&btsid=0ab6f83915858255934824697e4e73&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H9").Value  = "https://www.aliexpress.com/item/32847115923.html?spm=a2g0o.productlist.0.0.7bc147551VNoAO&algo_pvid=8df103e4-73f2-4f8f-817b-f8f95c907efc&algo_expid=8df103e4-73f2-4f8f-817b-f8f95c907efc-1&btsid=0ab6f82415858257023173189e2d5d&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H10").Value = "https://www.aliexpress.com/item/32847115923.html?spm=a2g0o.productlist.0.0.7bc147551VNoAO&algo_pvid=8df103e4-73f2-4f8f-817b-f8f95c907efc&algo_expid=8df103e4-73f2-4f8f-817b-f8f95c907efc-1&btsid=0ab6f82415858257023173189e2d5d&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H11").Value = "https://www.aliexpress.com/item/32847115923.html?spm=a2g0o.productlist.0.0.7bc147551VNoAO&algo_pvid=8df103e4-73f2-4f8f-817b-f8f95c907efc&algo_expid=8df103e4-73f2-4f8f-817b-f8f95c907efc-1&btsid=0ab6f82415858257023173189e2d5d&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H12").Value = "https://www.aliexpress.com/item/32369344670.html?spm=a2g0o.productlist.0.0.33481c2e3N4H2b&algo_pvid=88cedc09-4838-4e6a-8fbb-b9fe6ef80c69&algo_expid=88cedc09-4838-4e6a-8fbb-b9fe6ef80c69-4&btsid=0be3764315858258286375215ec63e&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"
$ws.Range("H14").Value = "https://www.aliexpress.com/item/32856564094.html?spm=a2g0o.productlist.0.0.74d42da5QBNTaI&algo_pvid=dc6c1b2d-0d11-482c-ab63-bf0ac9c0f492&algo_expid=dc6c1b2d-0d11-482c-ab63-bf0ac9c0f492-1&btsid=0ab6f81615858260205077678e29b5&ws_ab_test=searchweb0_0,searchweb201602_,searchweb201603_"

# ---------------------------------------------------------------------------
# 3. Mark the rows that are already confirmed/available (LDO1, LDO2,
#    ATTINY3216) with the built-in "Good" (green) cell style in the Mouser
#    column. Done before the hyperlink inserts below so the "Good" style
#    lands in the style table right after the pre-existing styles.
# ---------------------------------------------------------------------------
$ws.Range("G3").Style = "Good"
$ws.Range("G4").Style = "Good"
$ws.Range("G5").Style = "Good"

# ---------------------------------------------------------------------------
# 4. Add real hyperlinks on the Mouser (G) column for the three rows that
#    previously only had plain text (the URL text itself stays the same;
#    G2 and G3 already had hyperlinks and are left untouched).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G6"),  "https://www.mouser.be/ProductDetail/647-UWT1V101MCL1S")
$ws.Hyperlinks.Add($ws.Range("G7"),  "https://www.mouser.be/ProductDetail/963-JMK316ABJ107ML-T")
$ws.Hyperlinks.Add($ws.Range("G13"), "https://www.mouser.be/ProductDetail/581-SD1206T040S2R0")

# Re-assert the plain "Hyperlink" look (the .Add call above re-derives the
# font/style for the cell) so the three cells keep the same visual style as
# the pre-existing Mouser hyperlinks in G2/G3.
$ws.Range("G6").Style  = "Hyperlink"
$ws.Range("G7").Style  = "Hyperlink"
$ws.Range("G13").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 5. Column widths: column G got narrower (less text now that only the
#    Mouser link stays) and the new column H takes up the newly freed
#    space.
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 67.35
$ws.Columns.Item(8).ColumnWidth = 52.53

# Move the active selection the same way the author's session ended up.
$ws.Range("G19").Select()
